# Add one more day's worth of data (2020-05-14) to the 相談件数 sheet.
# This pushes the existing footer/note row (previously row 110) down to
# row 111, and fills the freed-up row 110 with the new day's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 110; Excel shifts the old row 110 (footer note) down
# to 111 and carries the formatting of row 109 onto the new row 110.
$ws.Rows.Item(110).Insert()

# Fill in the new day's data on row 110.
$ws.Cells.Item(110, 1).Value = 43965
$ws.Cells.Item(110, 2).Value = 229
$ws.Cells.Item(110, 3).Value = 37074
$ws.Cells.Item(110, 4).Value = 42
$ws.Cells.Item(110, 5).Value = 7530

# Keep the print area in sync with the newly extended data range.
$printArea = $wb.Names.Item(1)
$printArea.RefersTo = '=' + $ws.Name + '!$A$1:$E$112'

# Match the saved selection state (bottom-right pane now parks on A111).
$ws.Activate()
$ws.Range("A111").Select()
